$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2-10: 46065 -> 46066
foreach ($r in 2..10) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# Row 4: A 14516-2023 -> A 26262-2024
$ws.Cells.Item(4, 1).Value = "A 26262-2024"
$ws.Cells.Item(4, 2).Value = 45468.66077546297
$ws.Cells.Item(4, 7).Value = 0.6

# Row 5: A 26262-2024 -> A 14517-2023
$ws.Cells.Item(5, 1).Value = "A 14517-2023"
$ws.Cells.Item(5, 2).Value = 45012
$ws.Cells.Item(5, 7).Value = 0.6

# Row 6: A 4156-2023 -> A 14516-2023
$ws.Cells.Item(6, 1).Value = "A 14516-2023"
$ws.Cells.Item(6, 2).Value = 45012.86600694444
$ws.Cells.Item(6, 7).Value = 0.4

# Row 7: A 14517-2023 -> A 50762-2025
$ws.Cells.Item(7, 1).Value = "A 50762-2025"
$ws.Cells.Item(7, 2).Value = 45946
$ws.Cells.Item(7, 7).Value = 2.7

# Row 9: A 4159-2023 -> A 4156-2023
$ws.Cells.Item(9, 1).Value = "A 4156-2023"
$ws.Cells.Item(9, 2).Value = 44953
$ws.Cells.Item(9, 7).Value = 1.5

# Row 10: A 50762-2025 -> A 4159-2023
$ws.Cells.Item(10, 1).Value = "A 4159-2023"
$ws.Cells.Item(10, 2).Value = 44953
$ws.Cells.Item(10, 7).Value = 0.5
